$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04271373187048222
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.7270887267363939
